$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 290, shifting existing rows 290-409 down to 291-410
$ws.Rows.Item(290).Insert()

# Populate the new row 290 with the new record's data
$ws.Cells.Item(290, 1).Value = 5
$ws.Cells.Item(290, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(290, 3).Value = "Maule"
$ws.Cells.Item(290, 4).Value2 = 44837
$ws.Cells.Item(290, 5).Value = 7
$ws.Cells.Item(290, 6).Value = 100114013
$ws.Cells.Item(290, 7).Value = "Zanahoria"
$ws.Cells.Item(290, 8).Value = "Sin especificar"
$ws.Cells.Item(290, 9).Value = "Primera"
$ws.Cells.Item(290, 10).Value = 400
$ws.Cells.Item(290, 11).Value = 12000
$ws.Cells.Item(290, 12).Value = 12000
$ws.Cells.Item(290, 13).Value = 12000
$ws.Cells.Item(290, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(290, 15).Value = "Región de Ñuble"
$ws.Cells.Item(290, 16).Value = 600
$ws.Cells.Item(290, 17).Value = 20
$ws.Cells.Item(290, 18).Value = "Hortaliza"
